# "added missing reports logo"
#
# Slide 9 ("Tip: Welcome to When2Work") has an invisible click-rectangle
# ("Rectangle 7") drawn over the top "Start" button of the screenshot. It
# currently just advances to the next slide. We repoint it at the deck's
# first slide (the same target other "Rectangle"/"Rounded Rectangle" hot
# spots on the Error/Complete slides use) and add a second, identical hot
# spot ("Rectangle 8") over the bottom "Start" button that was missing it.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)

# Locate the existing hot-spot rectangle and the slide it should jump to.
$rect7 = $s.Shapes.Item("Rectangle 7")
$home = $p.Slides.Item(1)
$subAddress = $home.SlideID.ToString() + ",1," + $home.SlideIndex.ToString()

# Re-point Rectangle 7's click action at slide 1 instead of "next slide".
$link7 = $rect7.ActionSettings.Item(1).Hyperlink
$link7.SubAddress = $subAddress
$link7.Address = ""

# Duplicate it to create the missing second hot spot over the bottom
# "Start" button, then move/resize it into place.
$rect8 = $rect7.Duplicate()
$rect8.Name = "Rectangle 8"
$rect8.Left = 322.0832367464567
$rect8.Top = 416.15095528188976
$rect8.Width = 77.00866321732285
$rect8.Height = 19.903307986614173

# The duplicate already carries Rectangle 7's (now-updated) hyperlink, but
# set it explicitly too so both shapes unambiguously share the same target.
$link8 = $rect8.ActionSettings.Item(1).Hyperlink
$link8.SubAddress = $subAddress
$link8.Address = ""
